{"js": "// Remove the unused \"Abstract Title\" and \"Footnote Block Text\" paragraph\n// styles, and bump the \"Abstract\" style's space-before from 100 twips\n// (5pt) to 300 twips (15pt) so it matches its space-after.\nconst styles = context.document.getStyles();\n\n// 1) Delete the \"Abstract Title\" style.\nconst abstractTitle = styles.getByNameOrNullObject(\"Abstract Title\");\nabstractTitle.load(\"isNullObject\");\nawait context.sync();\nif (!abstractTitle.isNullObject) {\n  abstractTitle.delete();\n  await context.sync();\n}\n\n// 2) Delete the \"Footnote Block Text\" style.\nconst footnoteBlockText = styles.getByNameOrNullObject(\"Footnote Block Text\");\nfootnoteBlockText.load(\"isNullObject\");\nawait context.sync();\nif (!footnoteBlockText.isNullObject) {\n  footnoteBlockText.delete();\n  await context.sync();\n}\n\n// 3) Change the \"Abstract\" style's space-before to 300 twips (15pt).\nconst abstract = styles.getByNameOrNullObject(\"Abstract\");\nabstract.load(\"isNullObject\");\nawait context.sync();\nif (!abstract.isNullObject) {\n  abstract.paragraphFormat.spaceBefore = 15;\n  await context.sync();\n}\n", "ps1": "# Remove the unused \"Abstract Title\" and \"Footnote Block Text\" paragraph\n# styles, and bump the \"Abstract\" style's space-before from 100 twips\n# (5pt) to 300 twips (15pt) so it matches its space-after.\n$d = $word.ActiveDocument\n\n# 1) Delete the \"Abstract Title\" style.\ntry {\n    $abstractTitle = $d.Styles(\"Abstract Title\")\n    $abstractTitle.Delete()\n} catch {\n}\n\n# 2) Delete the \"Footnote Block Text\" style.\ntry {\n    $footnoteBlockText = $d.Styles(\"Footnote Block Text\")\n    $footnoteBlockText.Delete()\n} catch {\n}\n\n# 3) Change the \"Abstract\" style's space-before to 300 twips (15pt).\ntry {\n    $abstract = $d.Styles(\"Abstract\")\n    $abstract.ParagraphFormat.SpaceBefore = 15\n} catch {\n}\n"}
